$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 17 (fixture list now has 15 data rows instead of 16)
$ws.Rows.Item(17).Delete()

# Row 2
$ws.Range("A2").Value = 'Real Madrid  - Valencia CF: 19:00'
$ws.Range("B2").Value = 'Real Madrid'
$ws.Range("C2").Value = 78
$ws.Range("D2").Value = 87
$ws.Range("E2").Value = 92
$ws.Range("F2").Value = 1.18
$ws.Range("G2").Value = ""

# Row 3
$ws.Range("A3").Value = 'Bayern Munich  - Bayer 04 Leverkusen: 16:30'
$ws.Range("B3").Value = 'Bayern Munich'
$ws.Range("C3").Value = 74
$ws.Range("D3").Value = 88
$ws.Range("E3").Value = 83
$ws.Range("F3").Value = 1.23
$ws.Range("G3").Value = ""

# Row 4
$ws.Range("A4").Value = 'Burnley FC - Arsenal FC : -:-'''
$ws.Range("B4").Value = 'Arsenal FC'
$ws.Range("C4").Value = 71
$ws.Range("D4").Value = 74
$ws.Range("E4").Value = 83
$ws.Range("F4").Value = 1.26
$ws.Range("G4").Value = ""

# Row 5
$ws.Range("A5").Value = 'Haverfordwest County - The New Saints : -:-'''
$ws.Range("B5").Value = 'The New Saints'
$ws.Range("C5").Value = 71
$ws.Range("D5").Value = 78
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 1.67
$ws.Range("G5").Value = ""

# Row 6
$ws.Range("A6").Value = 'Lee Man ✓ - Hong Kong Football Club: 7:1'
$ws.Range("B6").Value = 'Lee Man'
$ws.Range("C6").Value = 60
$ws.Range("D6").Value = 72
$ws.Range("E6").Value = 68
$ws.Range("F6").Value = 2.05
$ws.Range("G6").Value = '✓'

# Row 7
$ws.Range("A7").Value = 'Atlético de Madrid  - Sevilla FC: 14:15'
$ws.Range("B7").Value = 'Atlético de Madrid'
$ws.Range("C7").Value = 54
$ws.Range("D7").Value = 87
$ws.Range("E7").Value = 92
$ws.Range("F7").Value = 1.18
$ws.Range("G7").Value = ""

# Row 8
$ws.Range("A8").Value = 'Liverpool FC Montevideo  - CA Juventud: 18:30'
$ws.Range("B8").Value = 'Liverpool FC Montevideo'
$ws.Range("C8").Value = 54
$ws.Range("D8").Value = 55
$ws.Range("E8").Value = 69
$ws.Range("F8").Value = 1.65
$ws.Range("G8").Value = ""

# Row 9
$ws.Range("A9").Value = 'Atlético Pantoja  - Atlántico FC: 21:00'
$ws.Range("B9").Value = 'Atlético Pantoja'
$ws.Range("C9").Value = 43
$ws.Range("D9").Value = 86
$ws.Range("E9").Value = 82
$ws.Range("F9").Value = 1.35
$ws.Range("G9").Value = ""

# Row 10
$ws.Range("A10").Value = 'Colwyn Bay  - Barry Town United: 16:15'
$ws.Range("B10").Value = 'Colwyn Bay'
$ws.Range("C10").Value = 42
$ws.Range("D10").Value = 88
$ws.Range("E10").Value = 83
$ws.Range("F10").Value = 1.23
$ws.Range("G10").Value = ""

# Row 11
$ws.Range("A11").Value = 'SSC Napoli  - Como 1907: 16:00'
$ws.Range("B11").Value = 'SSC Napoli'
$ws.Range("C11").Value = 42
$ws.Range("D11").Value = 87
$ws.Range("E11").Value = 76
$ws.Range("F11").Value = 1.91
$ws.Range("G11").Value = ""

# Row 12
$ws.Range("A12").Value = 'CD Real Cartagena  - Club Boca Juniors de Cali: 20:05'
$ws.Range("B12").Value = 'CD Real Cartagena'
$ws.Range("C12").Value = 41
$ws.Range("D12").Value = 87
$ws.Range("E12").Value = 92
$ws.Range("F12").Value = 1.18
$ws.Range("G12").Value = ""

# Row 13
$ws.Range("A13").Value = 'RB Leipzig  - VfB Stuttgart: -:-'''
$ws.Range("B13").Value = 'RB Leipzig'
$ws.Range("C13").Value = 38
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = 77
$ws.Range("F13").Value = 2.1
$ws.Range("G13").Value = ""

# Row 14
$ws.Range("A14").Value = 'Kryvbas Kryvyi Rig  - SC Poltava: 1:2'''
$ws.Range("B14").Value = 'Kryvbas Kryvyi Rig'
$ws.Range("C14").Value = 32
$ws.Range("D14").Value = 79
$ws.Range("E14").Value = 78
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = ""

# Row 15
$ws.Range("A15").Value = 'Ajax Amsterdam  - SC Heerenveen: 14:30'
$ws.Range("B15").Value = 'Ajax Amsterdam'
$ws.Range("C15").Value = 29
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = 90
$ws.Range("F15").Value = 1.55
$ws.Range("G15").Value = ""

# Row 16
$ws.Range("A16").Value = 'Real Oruro  - Gualberto Villarroel San José: 18:00'
$ws.Range("B16").Value = 'Real Oruro'
$ws.Range("C16").Value = 28
$ws.Range("D16").Value = 87
$ws.Range("E16").Value = 92
$ws.Range("F16").Value = 1.18
$ws.Range("G16").Value = ""
